$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("D3").Value = 44306
$ws.Range("M3").Value = 80

# Row 4
$ws.Range("D4").Value = 44313
$ws.Range("M4").Value = 120

# Row 6
$ws.Range("D6").Value = 44316
$ws.Range("M6").Value = 120

# Row 9
$ws.Range("D9").Value = 44330
$ws.Range("M9").Value = 60

# Row 10
$ws.Range("D10").Value = 44302
$ws.Range("M10").Value = 80
